$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: some "Price" values are plain decimals (e.g. "603.63") which Excel would
# otherwise auto-convert to a Number; a leading apostrophe (the doubled '' below)
# forces them to stay as Text, matching the source data's text cell type.
$ws.Range("D2").Value = '66.229.39'
$ws.Range("E2").Value = '  -0.07%  '

$ws.Range("D3").Value = '3.555.53'
$ws.Range("E3").Value = '  -0.23%  '

$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").Value = '''603.63'
$ws.Range("E5").Value = '  -0.43%  '

$ws.Range("D6").Value = '''146.99'
$ws.Range("E6").Value = '  +1.67%  '

$ws.Range("D7").Value = '3.557.10'
$ws.Range("E7").Value = '  -0.17%  '

$ws.Range("E8").Value = '  -0.22%  '

$ws.Range("E9").Value = '  -0.37%  '

$ws.Range("E10").Value = '  -1.70%  '

$ws.Range("D11").Value = '''7.87'
$ws.Range("E11").Value = '  +0.87%  '

$ws.Range("E12").Value = '  -0.90%  '

$ws.Range("D13").Value = '4.157.33'
$ws.Range("E13").Value = '  -0.31%  '

$ws.Range("D14").Value = '''0.0000202'
$ws.Range("E14").Value = '  -2.09%  '

$ws.Range("D15").Value = '''29.23'
$ws.Range("E15").Value = '  -3.79%  '

$ws.Range("D16").Value = '3.551.63'
$ws.Range("E16").Value = '  -0.20%  '

$ws.Range("E17").Value = '  +2.02%  '

$ws.Range("D18").Value = '66.201.35'
$ws.Range("E18").Value = '  -0.19%  '

$ws.Range("D19").Value = '''11.09'
$ws.Range("E19").Value = '  -3.52%  '

$ws.Range("E20").Value = '  +0.83%  '

$ws.Range("E21").Value = '  -0.25%  '

$ws.Range("D22").Value = '''419.43'
$ws.Range("E22").Value = '  -2.73%  '

$ws.Range("E23").Value = '  -0.97%  '

$ws.Range("D24").Value = '''77.81'
$ws.Range("E24").Value = '  -2.22%  '

$ws.Range("D25").Value = '3.693.81'
$ws.Range("E25").Value = '  -0.34%  '

$ws.Range("E26").Value = '  -0.06%  '

$ws.Range("E27").Value = '  -1.37%  '

$ws.Range("D28").Value = '''9.21'
$ws.Range("E28").Value = '  +0.60%  '

$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").Value = '''2.48'
$ws.Range("E29").Value = '  -1.05%  '

$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").Value = '''7.92'
$ws.Range("E30").Value = '  -0.54%  '

$ws.Range("D32").Value = '3.550.47'
$ws.Range("E32").Value = '  -0.26%  '

$ws.Range("E33").Value = '  +2.18%  '

$ws.Range("E34").Value = '  -3.02%  '

$ws.Range("E35").Value = '  -0.01%  '

$ws.Range("E36").Value = '  -2.43%  '

$ws.Range("E37").Value = '  -8.72%  '

$ws.Range("D38").Value = '''5.34'
$ws.Range("E38").Value = '  -4.77%  '

$ws.Range("E39").Value = '  -6.44%  '

$ws.Range("D40").Value = '''174.38'
$ws.Range("E40").Value = '  -0.78%  '

$ws.Range("D41").Value = '''0.0829'
$ws.Range("E41").Value = '  -2.34%  '

$ws.Range("E42").Value = '  -1.37%  '

$ws.Range("D43").Value = '''0.864'
$ws.Range("E43").Value = '  -2.86%  '

$ws.Range("D44").Value = '''45.76'
$ws.Range("E44").Value = '  -0.52%  '

$ws.Range("E45").Value = '  -4.97%  '

$ws.Range("E46").Value = '  +0.04%  '

$ws.Range("D47").Value = '''2.46'
$ws.Range("E47").Value = '  -2.99%  '

$ws.Range("D48").Value = '''7.15'
$ws.Range("E48").Value = '  -0.01%  '

$ws.Range("D49").Value = '''23.20'
$ws.Range("E49").Value = '  -1.03%  '

$ws.Range("E50").Value = '  -6.40%  '

$ws.Range("D51").Value = '''23.55'
$ws.Range("E51").Value = '  -6.30%  '
